$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(3)

# Bump the "Tuition Fees: ₹1,60,585" run sizes from 28pt to 32pt
# (only the first paragraph's two runs - leave the trailing blank
# paragraph alone).
$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$para1.Font.Size = 32

# Grow the text box to match the larger text (cy 800219 -> 861774 EMU).
$sh.Height = 67.85622047244094
